$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15
$ws.Cells.Item($row, 1).Value = 42620.888101851851
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item($row, 2).Value = 24
$ws.Cells.Item($row, 3).Value = 62
$ws.Cells.Item($row, 4).Value = 32
$ws.Cells.Item($row, 5).Value = 62
$ws.Cells.Item($row, 6).Value = 19
$ws.Cells.Item($row, 7).Value = 37727
$ws.Cells.Item($row, 8).Value = 18343
$ws.Cells.Item($row, 9).Value = 3251
$ws.Cells.Item($row, 10).Value = 393
$ws.Cells.Item($row, 11).Value = 200
$ws.Cells.Item($row, 12).Value = 37
$ws.Cells.Item($row, 13).Value = 9
$ws.Cells.Item($row, 14).Value = "Noun"
